# Weekly update: two new "Ají" (chili pepper) price records were reported
# for Terminal La Palmera de La Serena. They get inserted at the top of the
# existing block of "Ají" rows (rows 241-242), pushing the rest of the
# block down by two rows (old 241-262 -> new 243-264).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 241 so every following
# row (old 241..262) shifts down to (243..264), matching the target
# dimension A1:R264.
$ws.Rows("241:242").Insert()

# --- New row 241 -----------------------------------------------------
$ws.Range("A241").Value = 8
$ws.Range("B241").Value = "Terminal La Palmera de La Serena"
$ws.Range("C241").Value = "Coquimbo"
$ws.Range("D241").Value = 44783
$ws.Range("E241").Value = 4
$ws.Range("F241").Value = 100112021
$ws.Range("G241").Value = "Ají"
$ws.Range("H241").Value = "Inferno"
$ws.Range("I241").Value = "Primera"
$ws.Range("J241").Value = 460
$ws.Range("K241").Value = 14000
$ws.Range("L241").Value = 15000
$ws.Range("M241").Value = 14500
$ws.Range("N241").Value = "$/caja 12 kilos"
$ws.Range("O241").Value = "Región de Arica y Parinacota"
$ws.Range("P241").Value = 1208
$ws.Range("Q241").Value = 12
$ws.Range("R241").Value = "Hortaliza"

# --- New row 242 -----------------------------------------------------
$ws.Range("A242").Value = 8
$ws.Range("B242").Value = "Terminal La Palmera de La Serena"
$ws.Range("C242").Value = "Coquimbo"
$ws.Range("D242").Value = 44783
$ws.Range("E242").Value = 4
$ws.Range("F242").Value = 100112021
$ws.Range("G242").Value = "Ají"
$ws.Range("H242").Value = "Inferno"
$ws.Range("I242").Value = "Segunda"
$ws.Range("J242").Value = 280
$ws.Range("K242").Value = 8000
$ws.Range("L242").Value = 9000
$ws.Range("M242").Value = 8500
$ws.Range("N242").Value = "$/caja 12 kilos"
$ws.Range("O242").Value = "Región de Arica y Parinacota"
$ws.Range("P242").Value = 708
$ws.Range("Q242").Value = 12
$ws.Range("R242").Value = "Hortaliza"
